$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "subgenus" column (G) values to the new, more specific
# taxonomic order names ("frameshifting features" added host metadata).
$ws.Range("G2").Value = "Cetartiodactyla"
$ws.Range("G3").Value = "Primates"
$ws.Range("G4").Value = "Primates"
$ws.Range("G5").Value = "Primates"
$ws.Range("G6").Value = "Primates"
$ws.Range("G7").Value = "Chiroptera"

# G6 picks up the formatting used by the rest of column G (copy format
# from G2, whose style differs from the row-6 banding style).
$ws.Range("G2").Copy()
$ws.Range("G6").PasteSpecial(-4122)

# Move the active selection to A7.
$ws.Activate()
$ws.Range("A7").Select()
